# "added 4wk low sales check"
# Recalculates MyForecast (col D), Inventory Coverage (col H), Stockout Risk (col I),
# Reorder Urgency (col J) and Seasonality Index (col L) on the "Forecast Comparison"
# sheet, then refreshes the dependent roll-up numbers on the "Summary" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# ---- Forecast Comparison sheet -------------------------------------------------

# MyForecast (D2:D17) -> flattened to a constant 13 units/week
$ws1.Range("D2").Value  = 13
$ws1.Range("D3").Value  = 13
$ws1.Range("D4").Value  = 13
$ws1.Range("D5").Value  = 13
$ws1.Range("D6").Value  = 13
$ws1.Range("D7").Value  = 13
$ws1.Range("D8").Value  = 13
$ws1.Range("D9").Value  = 13
$ws1.Range("D11").Value = 13
$ws1.Range("D12").Value = 13
$ws1.Range("D13").Value = 13
$ws1.Range("D14").Value = 13
$ws1.Range("D15").Value = 13
$ws1.Range("D16").Value = 13
# D10 and D17 were already 13 - no change needed

# Inventory Coverage (H2:H9) - recomputed against the new forecast
$ws1.Range("H2").Value = 6.15
$ws1.Range("H3").Value = 5.15
$ws1.Range("H4").Value = 4.15
$ws1.Range("H5").Value = 3.15
$ws1.Range("H6").Value = 2.15
$ws1.Range("H7").Value = 1.15
$ws1.Range("H8").Value = 0.15
$ws1.Range("H9").Value = 0

# Stockout Risk (I8:I9) flips to High once coverage drops under the new 4wk low check
$ws1.Range("I8").Value = "High"
$ws1.Range("I9").Value = "High"

# Reorder Urgency (J8) flips to Urgent alongside the Stockout Risk change
$ws1.Range("J8").Value = "Urgent"

# Seasonality Index (L2:L17) - recomputed against the new forecast
$ws1.Range("L2").Value  = 1.13
$ws1.Range("L3").Value  = 0.81
$ws1.Range("L4").Value  = 0.86
$ws1.Range("L5").Value  = 0.93
$ws1.Range("L6").Value  = 1.06
$ws1.Range("L7").Value  = 0.88
$ws1.Range("L8").Value  = 1.07
$ws1.Range("L9").Value  = 0.97
$ws1.Range("L10").Value = 0.9
$ws1.Range("L11").Value = 0.82
$ws1.Range("L12").Value = 0.83
$ws1.Range("L13").Value = 0.95
$ws1.Range("L14").Value = 0.99
$ws1.Range("L15").Value = 0.96
$ws1.Range("L16").Value = 0.8
$ws1.Range("L17").Value = 0.84

# ---- Summary sheet --------------------------------------------------------------

$ws2.Range("B9").Value  = "208"   # Total Forecast (16 Weeks)
$ws2.Range("B10").Value = "104"   # Total Forecast (8 Weeks)
$ws2.Range("B11").Value = "52"    # Total Forecast (4 Weeks)
$ws2.Range("B12").Value = "13"    # Max Forecast
$ws2.Range("B14").Value = "13"    # Min Forecast

Write-Host "4wk low sales check applied"
